# Correção dos gráficos do subtema "Energia" do tema "Meio Ambiente"
# Convert the "Consumo de energia elétrica" column from Kwh to Mwh:
#   - header text: "(Kwh)" -> "(Mwh)"
#   - values: divide by 1000
#   - formatting: drop the special "#,##0 + shaded fill" number style that
#     used to visually group the big Kwh numbers, falling back to the plain
#     style already used by column A (no fill, General number format)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Header text: Kwh -> Mwh
$ws.Range("B1").Value = "Consumo de energia elétrica do município (Mwh)"

# 2) Re-format B2:B15 to match the plain style already used in column A
#    (no shaded fill, no #,##0 number format) by copying A2's formatting.
$ws.Range("A2").Copy()
$ws.Range("B2:B15").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# 3) Convert the raw Kwh figures to Mwh (divide by 1000)
$dataRange = $ws.Range("B2:B15")
for ($r = 1; $r -le $dataRange.Rows.Count; $r++) {
    $cell = $dataRange.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 / 1000
}
